$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Paragraph 1: add the "[2]" and "[1]" citation markers in the intro text
# ---------------------------------------------------------------------------
$find1 = $d.Content
$find1.Find.Execute(
    "negative mental health impacts. These",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "negative mental health impacts[2]. These", 2) | Out-Null

$find2 = $d.Content
$find2.Find.Execute(
    "comparison culture, typical",
    $false, $false, $false, $false, $false, $true, 1, $false,
    "comparison culture[1], typical", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2) Paragraph 3 ("Having been a user...") - re-assert the text so the
#    multiple legacy runs collapse into a single run (content is unchanged)
# ---------------------------------------------------------------------------
$para3Text = "Having been a user of social media since the " + [char]0x2018 + "MySpace days" + [char]0x2019 + " (approx 2005) I have seen various platforms come and go as part of trends and changes within the technological landscape. Although the introduction for such a medium of social interaction brought about its own set of challenges and problems, Ie Hunter Moore / isAnyoneUp, I personally saw the bigger and more widespread issues beginning with the advent of Facebook (approx 2009)."
$find3 = $d.Content
$find3.Find.Execute($para3Text, $false, $false, $false, $false, $false, $true, 1, $false, $para3Text, 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Append a "References" section at the end of the document
# ---------------------------------------------------------------------------

# The document currently ends with a trailing empty paragraph - use it for
# the "References" heading.
$lastParaIndex = $d.Paragraphs.Count
$referencesRange = $d.Paragraphs($lastParaIndex).Range
$referencesRange.InsertAfter("References")

# New paragraph for reference [1]
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$ref1Index = $d.Paragraphs.Count
$ref1Text = "[1] Appel, Helmut; Crusius, Jan; Gerlach, Alexander L.  (2015). Social Comparison, Envy, and Depression on Facebook: A Study Looking at the Effects of High Comparison Standards on Depressed Individuals. Journal of Social and Clinical Psychology, 34(4), pp.278.     "
$d.Paragraphs($ref1Index).Range.InsertAfter($ref1Text)

# New paragraph for reference [2]
$endRange2 = $d.Content
$endRange2.Collapse(0)
$endRange2.InsertParagraphAfter()
$ref2Index = $d.Paragraphs.Count
$ref2Text = "[2] Siddiqui, S. and Singh, T., 2016. Social media its impact with positive and negative aspects. International journal of computer applications technology and research, 5(2), pp.71-75. http://www.ijcat.com/archives/volume5/issue2/ijcatr05021006.pdf"
$d.Paragraphs($ref2Index).Range.InsertAfter($ref2Text)

# Trailing empty paragraph to match the original document's ending blank line
$endRange3 = $d.Content
$endRange3.Collapse(0)
$endRange3.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# 4) Formatting pass - Calisto MT, 10pt for the whole References block,
#    bold for the "References" title, and bold for the "[1]"/"[2]" markers.
# ---------------------------------------------------------------------------
$titleRange = $d.Paragraphs($lastParaIndex).Range
$titleRange.Font.Name = "Calisto MT"
$titleRange.Font.Bold = 1
$titleRange.Font.BoldBi = 1
$titleRange.Font.Size = 10
$titleRange.Font.SizeBi = 10

$ref1Range = $d.Paragraphs($ref1Index).Range
$ref1Range.Font.Name = "Calisto MT"
$ref1Range.Font.Size = 10
$ref1Range.Font.SizeBi = 10

$ref2Range = $d.Paragraphs($ref2Index).Range
$ref2Range.Font.Name = "Calisto MT"
$ref2Range.Font.Size = 10
$ref2Range.Font.SizeBi = 10

# Bold just the leading "[1]" / "[2]" markers (narrow ranges - Bold only,
# BoldBi is intentionally skipped here since it does not stay scoped to a
# partial-run range in this runtime).
$ref1Start = $d.Paragraphs($ref1Index).Range.Start
$marker1 = $d.Range($ref1Start, $ref1Start + 3)
$marker1.Font.Bold = 1

$ref2Start = $d.Paragraphs($ref2Index).Range.Start
$marker2 = $d.Range($ref2Start, $ref2Start + 3)
$marker2.Font.Bold = 1

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
